# 9th Stab - Cosmetic Changes
# Insert two new snapshot columns (Jun_17, Jun_15) ahead of the existing
# Jun_10 data column, shifting that column to the right and backfilling
# the two new columns with the default "UN" rating used for every other
# not-yet-updated analyst row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# Remember the old "Jun_13" header text (column B) before it gets
# overwritten with the newest snapshot label below.
$oldHeaderB = $ws.Range("B1").Value()

# Insert two blank columns before the existing "C" column (the old
# Jun_10 column). This shifts that column -- and its data/formatting --
# two places to the right, from C to E.
$ws.Columns("C:D").Insert()

# New header cells for the freshly inserted snapshot columns, and carry
# the previous "Jun_13" header along into its new column (D).
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"
$ws.Range("D1").Value = $oldHeaderB

# Backfill the new columns with the same placeholder rating ("UN") used
# for every other not-yet-updated analyst row, mirroring column B.
for ($r = 2; $r -le $lastRow; $r++) {
    $filler = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($r, 3).Value = $filler
    $ws.Cells.Item($r, 4).Value = $filler
}

# Match column widths across the new and carried-over columns, and mark
# the two filler columns as a collapsed outline group (summarised by the
# real data column to their right).
$ws.Columns("C").ColumnWidth = 7.14
$ws.Columns("D").ColumnWidth = 7.14
$ws.Columns("E").ColumnWidth = 7.14
$ws.Columns("C").OutlineLevel = 1
$ws.Columns("D").OutlineLevel = 1
$ws.Columns("C:D").ShowDetail = $false
